$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(42613.760983796295, 18, 54, 40, 54, 23, 29384, 18956, 3136, 418, 309, 30, 9, "Noun"),
    @(42613.8906712963,   20, 54, 42, 54, 25, 14162, 17650, 2892, 392, 306, 35, 12, "Noun"),
    @(42614.887499999997, 40, 63, 35, 63, 10, 22246, 12256, 2006, 289, 160, 42, 5,  "Noun"),
    @(42615.886770833335, 32, 60, 34, 60, 9,  14718, 15505, 2608, 375, 216, 39, 4,  "Noun")
)

$startRow = 9
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 0; $col -lt $values.Count; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $values[$col]
    }
}
